$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data last refreshed" timestamp shown in row 1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 11:52"

# Espana (row 5) - refreshed case counts
$ws.Range("B5").Value = 169496
$ws.Range("C5").Value = 2665
$ws.Range("D5").Value = 64727
$ws.Range("E5").Value = 87280
$ws.Range("G5").Value = 280
$ws.Range("H5").Value = 17489

# Malasia (row 38) - refreshed case counts
$ws.Range("B38").Value = 4817
$ws.Range("C38").Value = 134
$ws.Range("D38").Value = 2276
$ws.Range("E38").Value = 2464
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 77

# Kuwait (row 67) - refreshed case counts
$ws.Range("B67").Value = 1300
$ws.Range("C67").Value = 66
$ws.Range("D67").Value = 150
$ws.Range("E67").Value = 1148
$ws.Range("F67").Value = 26
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 2

# Libano (row 89) - refreshed case counts
$ws.Range("B89").Value = 632
$ws.Range("C89").Value = 2
$ws.Range("E89").Value = 532

# Rows 139/140: Jamaica and Etiopia swap order (Etiopia now listed first)
# and both get refreshed case counts
$ws.Range("A139").Value = "Etiopia"
$ws.Range("B139").Value = 74
$ws.Range("C139").Value = 3
$ws.Range("D139").Value = 14
$ws.Range("E139").Value = 57
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 3

$ws.Range("A140").Value = "Jamaica"
$ws.Range("B140").Value = 72
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 16
$ws.Range("E140").Value = 52
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 4

# Rows 196-200: Cabo Verde moves up ahead of Surinam, shifting
# Surinam / Islas Turcas y Caicos / Gambia / Nicaragua down by one row,
# all with refreshed case counts
$ws.Range("A196").Value = "Cabo Verde"
$ws.Range("B196").Value = 10
$ws.Range("C196").Value = 2
$ws.Range("D196").Value = 1
$ws.Range("E196").Value = 8
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 1

$ws.Range("A197").Value = "Surinam"
$ws.Range("B197").Value = 10
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 4
$ws.Range("E197").Value = 5
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 1

$ws.Range("A198").Value = "Islas Turcas y Caicos"
$ws.Range("B198").Value = 9
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 0
$ws.Range("E198").Value = 8
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Gambia"
$ws.Range("B199").Value = 9
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 2
$ws.Range("E199").Value = 6
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

$ws.Range("A200").Value = "Nicaragua"
$ws.Range("B200").Value = 9
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 4
$ws.Range("E200").Value = 4
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 1
